# Append a new row (61) of data to each of the 4 worksheets.
$wb = $excel.ActiveWorkbook

$rowsData = @{
    "ROW35-FE-LIFTER" = @{
        A = 45754.8737634375
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x6e"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 366
        I = 13
    }
    "ROW35-MID-LIFTER" = @{
        A = 45754.72441709491
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x6a"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 362
        I = 14
    }
    "ROW02-FE-LIFTER" = @{
        A = 45754.86115554398
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x6e"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 366
        I = 3
    }
    "ROW02-MID-LIFTER" = @{
        A = 45754.92579674769
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x6a"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 362
        I = 3
    }
}

foreach ($sheetName in $rowsData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = $rowsData[$sheetName]

    # find next empty row after the last used row in column A
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    $ws.Cells.Item($newRow, 1).Value = $row.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $row.B
    $ws.Cells.Item($newRow, 3).Value = $row.C
    $ws.Cells.Item($newRow, 4).Value = $row.D
    $ws.Cells.Item($newRow, 5).Value = $row.E
    $ws.Cells.Item($newRow, 6).Value = $row.F
    $ws.Cells.Item($newRow, 7).Value = $row.G
    $ws.Cells.Item($newRow, 8).Value = $row.H
    $ws.Cells.Item($newRow, 9).Value = $row.I
}
